{"js": "// Resume edit: \"Managed\" -> \"Led\" (site-wide + multi-client project bullets),\n// and the title \"Web/Software Architect\" -> \"Web Architect/Lead Developer\".\n//\n// Using Body.search + Range.insertText(..., Replace) keeps each matched\n// run's original character formatting (font, color, size) intact, matching\n// the author's intent of only changing the wording.\n\n// 1) Title line: \"Web/Software Architect\" -> \"Web Architect/Lead Developer\"\nconst title = context.document.body.search(\"Web/Software Architect\", {\n  matchCase: true,\n  matchWholeWord: false\n});\ntitle.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < title.items.length; i++) {\n  title.items[i].insertText(\"Web Architect/Lead Developer\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"Managed\" -> \"Led\" everywhere it appears as its own word\n//    (\"Managed site-wide performance improvement projects\" and\n//    \"Managed a multi-client checkout flow performance improvement project\")\nconst managed = context.document.body.search(\"Managed\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nmanaged.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < managed.items.length; i++) {\n  managed.items[i].insertText(\"Led\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Resume edit: \"Managed\" -> \"Led\" (site-wide + multi-client project bullets),\n# and the title \"Web/Software Architect\" -> \"Web Architect/Lead Developer\".\n#\n# Using Find/Replace (wdReplaceAll) keeps each matched run's original\n# character formatting (font, color, size) intact, matching the author's\n# intent of only changing the wording, not the style.\n\n$d = $word.ActiveDocument\n\n# 1) Title line: \"Web/Software Architect\" -> \"Web Architect/Lead Developer\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Web/Software Architect\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Web Architect/Lead Developer\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n\n# 2) \"Managed\" -> \"Led\" everywhere it appears as its own word\n#    (\"Managed site-wide performance improvement projects\" and\n#    \"Managed a multi-client checkout flow performance improvement project\")\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Managed\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Led\"\n$find2.Execute([ref]$find2.Text, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2)\n"}
